# Applies the 2023-12-24 cryptos-list refresh: updated prices/24h deltas for
# existing rows, plus a new "LEO" row inserted at row 28 which pushes the
# rest of the list down by one (the final row, WOONetwork, drops off the end).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = '''43.607.58'
$ws.Cells.Item(2, 5).Value = '  +0.04%  '

# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = '''2.289.47'
$ws.Cells.Item(3, 5).Value = '  +0.12%  '

# Row 4: TetherUSD
$ws.Cells.Item(4, 5).Value = '  +0.22%  '

# Row 5: Solana
$ws.Cells.Item(5, 4).Value = '''110.72'
$ws.Cells.Item(5, 5).Value = '  +15.35%  '

# Row 6: BNB
$ws.Cells.Item(6, 4).Value = '''267.54'
$ws.Cells.Item(6, 5).Value = '  -0.07%  '

# Row 7: XRP
$ws.Cells.Item(7, 5).Value = '  +0.45%  '

# Row 8: USDC
$ws.Cells.Item(8, 5).Value = '  +0.27%  '

# Row 9: Cardano
$ws.Cells.Item(9, 4).Value = '''0.614'
$ws.Cells.Item(9, 5).Value = '  +0.58%  '

# Row 10: Avalanche
$ws.Cells.Item(10, 4).Value = '''47.40'
$ws.Cells.Item(10, 5).Value = '  +4.11%  '

# Row 11: Dogecoin
$ws.Cells.Item(11, 4).Value = '''0.0945'
$ws.Cells.Item(11, 5).Value = '  +1.13%  '

# Row 12: Polkadot
$ws.Cells.Item(12, 4).Value = '''8.95'
$ws.Cells.Item(12, 5).Value = '  +12.14%  '

# Row 13: TRON
$ws.Cells.Item(13, 5).Value = '  +0.46%  '

# Row 14: Chainlink
$ws.Cells.Item(14, 4).Value = '''15.73'
$ws.Cells.Item(14, 5).Value = '  +2.89%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Cells.Item(15, 4).Value = '''2.629.39'
$ws.Cells.Item(15, 5).Value = '  -0.09%  '

# Row 16: Polygon
$ws.Cells.Item(16, 4).Value = '''0.844'
$ws.Cells.Item(16, 5).Value = '  -0.55%  '

# Row 17: WrappedEther
$ws.Cells.Item(17, 4).Value = '''2.284.03'
$ws.Cells.Item(17, 5).Value = '  -0.15%  '

# Row 18: WrappedBTC
$ws.Cells.Item(18, 4).Value = '''43.607.19'
$ws.Cells.Item(18, 5).Value = '  +0.01%  '

# Row 19: ShibaInu
$ws.Cells.Item(19, 5).Value = '  +0.36%  '

# Row 20: Uniswap
$ws.Cells.Item(20, 4).Value = '''6.74'
$ws.Cells.Item(20, 5).Value = '  +8.82%  '

# Row 21: Litecoin
$ws.Cells.Item(21, 4).Value = '''72.18'
$ws.Cells.Item(21, 5).Value = '  -0.12%  '

# Row 22: ImmutableX
$ws.Cells.Item(22, 4).Value = '''2.46'
$ws.Cells.Item(22, 5).Value = '  -4.31%  '

# Row 23: BitcoinCash
$ws.Cells.Item(23, 4).Value = '''232.04'
$ws.Cells.Item(23, 5).Value = '  -0.34%  '

# Row 24: InternetComputer(DFINITY)
$ws.Cells.Item(24, 4).Value = '''9.65'
$ws.Cells.Item(24, 5).Value = '  +6.27%  '

# Row 25: PancakeSwap
$ws.Cells.Item(25, 5).Value = '  +9.03%  '

# Row 26: Dai
$ws.Cells.Item(26, 5).Value = '  +0.05%  '

# Row 27: Cosmos
$ws.Cells.Item(27, 4).Value = '''11.60'
$ws.Cells.Item(27, 5).Value = '  +3.73%  '

# Row 28: LEO
$ws.Cells.Item(28, 2).Value = 'LEO'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(28, 4).Value = '''3.92'
$ws.Cells.Item(28, 5).Value = '  +0.99%  '

# Row 29: InjectiveProtocol
$ws.Cells.Item(29, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(29, 4).Value = '''41.73'
$ws.Cells.Item(29, 5).Value = '  +4.23%  '

# Row 30: WEMIXToken
$ws.Cells.Item(30, 2).Value = 'WEMIXToken'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(30, 4).Value = '''3.39'
$ws.Cells.Item(30, 5).Value = '  -2.20%  '

# Row 31: Toncoin
$ws.Cells.Item(31, 2).Value = 'Toncoin'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(31, 4).Value = '''2.26'
$ws.Cells.Item(31, 5).Value = '  -1.06%  '

# Row 32: Monero
$ws.Cells.Item(32, 2).Value = 'Monero'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(32, 4).Value = '''176.03'
$ws.Cells.Item(32, 5).Value = '  +0.53%  '

# Row 33: EthereumClassic
$ws.Cells.Item(33, 2).Value = 'EthereumClassic'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(33, 4).Value = '''21.54'
$ws.Cells.Item(33, 5).Value = '  -1.32%  '

# Row 34: Hedera
$ws.Cells.Item(34, 2).Value = 'Hedera'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(34, 4).Value = '''0.0925'
$ws.Cells.Item(34, 5).Value = '  +3.29%  '

# Row 35: Filecoin
$ws.Cells.Item(35, 2).Value = 'Filecoin'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(35, 4).Value = '''5.61'
$ws.Cells.Item(35, 5).Value = '  +4.86%  '

# Row 36: Stellar
$ws.Cells.Item(36, 2).Value = 'Stellar'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(36, 4).Value = '''0.127'
$ws.Cells.Item(36, 5).Value = '  +0.92%  '

# Row 37: RenderToken
$ws.Cells.Item(37, 2).Value = 'RenderToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(37, 4).Value = '''4.72'
$ws.Cells.Item(37, 5).Value = '  +7.80%  '

# Row 38: VeChain
$ws.Cells.Item(38, 2).Value = 'VeChain'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(38, 4).Value = '''0.0361'
$ws.Cells.Item(38, 5).Value = '  +2.23%  '

# Row 39: Kaspa
$ws.Cells.Item(39, 2).Value = 'Kaspa'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(39, 4).Value = '''0.107'
$ws.Cells.Item(39, 5).Value = '  -0.57%  '

# Row 40: NEARProtocol
$ws.Cells.Item(40, 2).Value = 'NEARProtocol'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(40, 4).Value = '''3.77'
$ws.Cells.Item(40, 5).Value = '  +12.73%  '

# Row 41: LidoDAOToken
$ws.Cells.Item(41, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(41, 4).Value = '''2.41'
$ws.Cells.Item(41, 5).Value = '  +5.02%  '

# Row 42: Algorand
$ws.Cells.Item(42, 2).Value = 'Algorand'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(42, 4).Value = '''0.242'
$ws.Cells.Item(42, 5).Value = '  +0.71%  '

# Row 43: Celestia
$ws.Cells.Item(43, 2).Value = 'Celestia'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Cells.Item(43, 4).Value = '''13.65'
$ws.Cells.Item(43, 5).Value = '  +10.98%  '

# Row 44: MultiversX
$ws.Cells.Item(44, 2).Value = 'MultiversX'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Cells.Item(44, 4).Value = '''72.29'
$ws.Cells.Item(44, 5).Value = '  +10.53%  '

# Row 45: THORChain
$ws.Cells.Item(45, 2).Value = 'THORChain'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Cells.Item(45, 4).Value = '''6.24'
$ws.Cells.Item(45, 5).Value = '  +20.60%  '

# Row 46: FirstDigitalUSD
$ws.Cells.Item(46, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(46, 4).Value = '''1.00'
$ws.Cells.Item(46, 5).Value = '  +0.04%  '

# Row 47: ARBITRUM
$ws.Cells.Item(47, 2).Value = 'ARBITRUM'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(47, 4).Value = '''1.38'
$ws.Cells.Item(47, 5).Value = '  +2.22%  '

# Row 48: FraxShare
$ws.Cells.Item(48, 2).Value = 'FraxShare'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(48, 4).Value = '''8.83'
$ws.Cells.Item(48, 5).Value = '  +0.55%  '

# Row 49: Aave
$ws.Cells.Item(49, 4).Value = '''102.25'
$ws.Cells.Item(49, 5).Value = '  +5.15%  '

# Row 50: Cronos
$ws.Cells.Item(50, 2).Value = 'Cronos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(50, 4).Value = '''0.0995'
$ws.Cells.Item(50, 5).Value = '  -2.09%  '

# Row 51: TrustWalletToken
$ws.Cells.Item(51, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(51, 4).Value = '''1.22'
$ws.Cells.Item(51, 5).Value = '  +2.26%  '
